$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.313328742980957
$ws.Range("B1").Value = 3.4656982421875
$ws.Range("C1").Value = 3.187313556671143
$ws.Range("D1").Value = 1.580045461654663
$ws.Range("E1").Value = 1.139871954917908
